# Generate Report for Handoff
# Updates the "Latest Handoff Date/Datetime" entries for the
# d4f4aaa7-593f-474a-9856-2802a19e9e93.md file (row 7) across the
# Overview, zh-cn and de-de sheets to reflect a new handoff event.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("D7").Value = "2016-41-18 05:41:29"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E7").Value = "2016-03-18 05:41:27"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E7").Value = "2016-03-18 05:41:29"
